# Adds a new "2022-Q3" worksheet (fund-holding detail) right after the
# "2022-Q2" sheet's position (i.e. as the 2nd sheet, after "总计"), and
# inserts a corresponding summary row at the top of the "总计" table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: assign literal TEXT to a cell without Excel's COM layer
# silently re-interpreting numeric-looking strings (e.g. "6.37", "005396")
# as numbers. We build the text via a throw-away formula cell (whose
# result type is always string), copy it, and paste-special *values only*
# into the destination - this keeps the destination's native text type.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($ws, $range, [string]$text)
    $helper = $ws.Range("ZZ1000")
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)
    $helper.ClearContents()
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the existing "2022-Q2"
#    sheet (same column layout/styles) and placing the copy immediately
#    before it - this reproduces the target's column formatting with no
#    manual style plumbing.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$ws = $wb.Worksheets.Item("2022-Q2 (2)")
$ws.Name = "2022-Q3"

# The template only has 3 data rows (rows 2-4); the new sheet needs 6
# (rows 2-7). Extend it by duplicating the row-2 pattern (style + column
# layout) down into rows 5-7.
$ws.Range("A2:H2").Copy($ws.Range("A5:H5"))
$ws.Range("A2:H2").Copy($ws.Range("A6:H6"))
$ws.Range("A2:H2").Copy($ws.Range("A7:H7"))

# Fund-holding rows for 2022-Q3.
$data = @(
    @(0, "159610", "景顺长城中证500增强策略ETF", "6.37", "98.41", "1.16", "0.0739", 8),
    @(1, "005396", "中金丰硕混合",                 "1.77", "76.61", "4.06", "0.0719", 9),
    @(2, "970042", "国海量化优选一年持有股票C",     "7.16", "87.31", "0.33", "0.0236", 10),
    @(3, "009613", "上银中证500指数增强A",         "0.97", "92.48", "1.00", "0.0097", 9),
    @(4, "009614", "上银中证500指数增强C",         "0.74", "92.48", "1.00", "0.0074", 9),
    @(5, "970041", "国海量化优选一年持有股票A",     "0.63", "87.31", "0.33", "0.0021", 10)
)

$r = 2
foreach ($row in $data) {
    $ws.Range("A$r").Value = $row[0]
    Set-TextValue $ws $ws.Range("B$r") $row[1]
    Set-TextValue $ws $ws.Range("C$r") $row[2]
    Set-TextValue $ws $ws.Range("D$r") $row[3]
    Set-TextValue $ws $ws.Range("E$r") $row[4]
    Set-TextValue $ws $ws.Range("F$r") $row[5]
    Set-TextValue $ws $ws.Range("G$r") $row[6]
    $ws.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert the new summary row into "总计": shift the existing rows 2-8
#    down to 3-9 (copying preserves per-cell style/type), then write the
#    new 2022-Q3 summary into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
for ($row = 8; $row -ge 2; $row--) {
    $src = $summary.Range("A" + $row + ":D" + $row)
    $dst = $summary.Range("A" + ($row + 1) + ":D" + ($row + 1))
    $src.Copy($dst)
}

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.19

# Re-number the index column (A) for the shifted rows so it stays a
# simple 0-based sequence.
$idx = 1
for ($row = 3; $row -le 9; $row++) {
    $summary.Range("A$row").Value = $idx
    $idx = $idx + 1
}
